$d = $word.ActiveDocument

$d.Content.Find.Execute("976÷2=488, 0", $true, $false, $false, $false, $false, $true, 1, $false, "745÷4=186, 1", 2) | Out-Null
$d.Content.Find.Execute("308÷4=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "786÷6=131, 0", 2) | Out-Null
$d.Content.Find.Execute("535÷7=76, 3", $true, $false, $false, $false, $false, $true, 1, $false, "979÷4=244, 3", 2) | Out-Null
$d.Content.Find.Execute("358÷8=44, 6", $true, $false, $false, $false, $false, $true, 1, $false, "594÷4=148, 2", 2) | Out-Null
$d.Content.Find.Execute("604÷7=86, 2", $true, $false, $false, $false, $false, $true, 1, $false, "805÷4=201, 1", 2) | Out-Null
$d.Content.Find.Execute("605÷2=302, 1", $true, $false, $false, $false, $false, $true, 1, $false, "848÷2=424, 0", 2) | Out-Null
$d.Content.Find.Execute("547÷4=136, 3", $true, $false, $false, $false, $false, $true, 1, $false, "273÷3=91, 0", 2) | Out-Null
$d.Content.Find.Execute("129÷9=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "173÷9=19, 2", 2) | Out-Null
$d.Content.Find.Execute("448÷2=224, 0", $true, $false, $false, $false, $false, $true, 1, $false, "540÷7=77, 1", 2) | Out-Null
$d.Content.Find.Execute("238÷9=26, 4", $true, $false, $false, $false, $false, $true, 1, $false, "783÷7=111, 6", 2) | Out-Null
$d.Content.Find.Execute("251÷5=50, 1", $true, $false, $false, $false, $false, $true, 1, $false, "361÷8=45, 1", 2) | Out-Null
$d.Content.Find.Execute("584÷8=73, 0", $true, $false, $false, $false, $false, $true, 1, $false, "952÷3=317, 1", 2) | Out-Null
$d.Content.Find.Execute("839÷8=104, 7", $true, $false, $false, $false, $false, $true, 1, $false, "506÷6=84, 2", 2) | Out-Null
$d.Content.Find.Execute("530÷2=265, 0", $true, $false, $false, $false, $false, $true, 1, $false, "927÷4=231, 3", 2) | Out-Null
$d.Content.Find.Execute("491÷8=61, 3", $true, $false, $false, $false, $false, $true, 1, $false, "527÷7=75, 2", 2) | Out-Null
$d.Content.Find.Execute("448÷8=56, 0", $true, $false, $false, $false, $false, $true, 1, $false, "892÷5=178, 2", 2) | Out-Null
$d.Content.Find.Execute("597÷8=74, 5", $true, $false, $false, $false, $false, $true, 1, $false, "799÷7=114, 1", 2) | Out-Null
$d.Content.Find.Execute("485÷5=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "956÷2=478, 0", 2) | Out-Null
$d.Content.Find.Execute("333÷8=41, 5", $true, $false, $false, $false, $false, $true, 1, $false, "241÷5=48, 1", 2) | Out-Null
$d.Content.Find.Execute("286÷2=143, 0", $true, $false, $false, $false, $false, $true, 1, $false, "235÷2=117, 1", 2) | Out-Null
$d.Content.Find.Execute("150÷9=16, 6", $true, $false, $false, $false, $false, $true, 1, $false, "383÷6=63, 5", 2) | Out-Null
$d.Content.Find.Execute("579÷8=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "131÷7=18, 5", 2) | Out-Null
$d.Content.Find.Execute("908÷7=129, 5", $true, $false, $false, $false, $false, $true, 1, $false, "165÷7=23, 4", 2) | Out-Null
$d.Content.Find.Execute("106÷7=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "104÷6=17, 2", 2) | Out-Null
$d.Content.Find.Execute("253÷7=36, 1", $true, $false, $false, $false, $false, $true, 1, $false, "809÷5=161, 4", 2) | Out-Null
